# Basic accel value passing
# Adds Sheet2 with accelerometer summary data (X/Y/Z peak + count), an
# autofilter, and per-column color-scale conditional formatting; mirrors
# the activeTab / tabSelected bookkeeping Excel performs when the new
# sheet becomes the active one.

$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item(1)

# Insert the new sheet right after Sheet1 so it becomes Worksheets(2) and
# (being freshly added/activated) the workbook's active tab.
$ws2 = $wb.Worksheets.Add($null, $ws1)
$ws2.Name = "Sheet2"

# --- Data rows (written first, so "Count: " is registered in the shared
#     string table before the "X"/"Y"/"Z" header labels below) ---
$data = @"
-352,4028,856,7
-332,4044,824,8
1294,3662,504,9
4138,1610,-12,10
-3976,2768,786,11
1820,3044,198,12
-2924,3094,-6,13
1192,3334,-862,14
-1732,2296,-172,15
1990,5228,68,16
-3232,3196,1762,17
-1466,2468,-120,18
974,3810,1704,19
-4014,4116,1052,20
1124,4760,-1270,21
-1120,3200,3342,22
-1310,5480,-674,23
-844,2174,-4346,24
-670,2518,1624,25
-1224,3306,790,26
528,2816,-1540,27
-1644,2332,1232,28
2062,2498,-224,29
-1682,3856,-1002,30
-2076,3128,2238,31
-2424,6114,3454,32
3372,2028,-820,33
-1466,3220,756,34
-814,2064,-812,35
5438,1260,-148,36
-3252,4134,-356,37
-1020,2618,3996,38
-376,2804,-1086,39
-544,2332,-1694,40
-2782,2350,2448,41
-206,7032,-1560,42
5080,624,1642,43
-2196,4724,-202,44
3292,4882,-90,45
-1204,5988,710,46
832,2794,-26,47
-776,2690,1650,48
-786,3890,-1636,49
"@ -split "`n"

$row = 2
foreach ($line in $data) {
    $parts = $line.Trim() -split ","
    $x = [double]$parts[0]
    $y = [double]$parts[1]
    $z = [double]$parts[2]
    $cnt = [double]$parts[3]

    $ws2.Cells.Item($row, 1).Value = "X: "
    $ws2.Cells.Item($row, 2).Value = $x
    $ws2.Cells.Item($row, 3).Value = "Y: "
    $ws2.Cells.Item($row, 4).Value = $y
    $ws2.Cells.Item($row, 5).Value = "Z: "
    $ws2.Cells.Item($row, 6).Value = $z
    $ws2.Cells.Item($row, 7).Value = "Count: "
    $ws2.Cells.Item($row, 8).Value = $cnt

    $row++
}

# --- Header row, written last on purpose (see note above) ---
$ws2.Cells.Item(1, 2).Value = "X"
$ws2.Cells.Item(1, 4).Value = "Y"
$ws2.Cells.Item(1, 6).Value = "Z"

# --- Selection matches the authored file (cell D8 active on Sheet2) ---
[void]$ws2.Range("D8").Select()

# --- AutoFilter over the full data range, with the hidden sheet-scoped
#     _FilterDatabase defined name Excel normally creates alongside it ---
$usedRange = $ws2.Range("A1:H44")
[void]$usedRange.AutoFilter()
$filterName = $ws2.Names.Add("_xlnm._FilterDatabase", "=Sheet2!`$A`$1:`$H`$44")
$filterName.Visible = $false

# --- Per-column 3-color scale conditional formatting (B, D, F; whole
#     column, matching the authored file's 1:1048576 sqref). Priorities
#     are set explicitly afterwards (3/2/1 for B/D/F) to match the
#     authored file. ---
$colB = $ws2.Range("B1:B1048576")
$cfB = $colB.FormatConditions.AddColorScale(3)
$cfB.ColorScaleCriteria.Item(1).Type = 1
$cfB.ColorScaleCriteria.Item(1).FormatColor.Color = 7039480
$cfB.ColorScaleCriteria.Item(2).Type = 4
$cfB.ColorScaleCriteria.Item(2).Value = 50
$cfB.ColorScaleCriteria.Item(2).FormatColor.Color = 8711167
$cfB.ColorScaleCriteria.Item(3).Type = 2
$cfB.ColorScaleCriteria.Item(3).FormatColor.Color = 8109667

$colD = $ws2.Range("D1:D1048576")
$cfD = $colD.FormatConditions.AddColorScale(3)
$cfD.ColorScaleCriteria.Item(1).Type = 1
$cfD.ColorScaleCriteria.Item(1).FormatColor.Color = 7039480
$cfD.ColorScaleCriteria.Item(2).Type = 4
$cfD.ColorScaleCriteria.Item(2).Value = 50
$cfD.ColorScaleCriteria.Item(2).FormatColor.Color = 8711167
$cfD.ColorScaleCriteria.Item(3).Type = 2
$cfD.ColorScaleCriteria.Item(3).FormatColor.Color = 8109667

$colF = $ws2.Range("F1:F1048576")
$cfF = $colF.FormatConditions.AddColorScale(3)
$cfF.ColorScaleCriteria.Item(1).Type = 1
$cfF.ColorScaleCriteria.Item(1).FormatColor.Color = 7039480
$cfF.ColorScaleCriteria.Item(2).Type = 4
$cfF.ColorScaleCriteria.Item(2).Value = 50
$cfF.ColorScaleCriteria.Item(2).FormatColor.Color = 8711167
$cfF.ColorScaleCriteria.Item(3).Type = 2
$cfF.ColorScaleCriteria.Item(3).FormatColor.Color = 8109667

$cfB.Priority = 3
$cfD.Priority = 2
$cfF.Priority = 1
